# Rename the 9 worksheets in place (positions stay the same, names change)
# as part of rerunning LU d2c FeatEng for FR cities with new spatial units.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ34534036",
    "summ34900327",
    "summ35250164",
    "summ35655179",
    "summ36002040",
    "summ36338418",
    "summ36671391",
    "summ37020463",
    "summ37365029"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
